$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RevA-B")

# Header row - new "Alternate PN" / "Alternate Digikey PN" columns
$ws.Range("I3").Value = "Alternate PN"
$ws.Range("J3").Value = "Alternate Digikey PN"

# Row 4 - MAX30101 sensor alternate: MAX30100EFD+ / MAX30100EFD+TCT-ND
$ws.Range("J4").Value = "MAX30100EFD+TCT-ND"

# Row 5 - Level shifter alternate: TCA9801DGKT / 296-46571-1-ND
$ws.Range("J5").Value = "296-46571-1-ND"
$ws.Range("I5").Value = "TCA9801DGKT"

# Row 6 - Voltage regulator alternate: NCP508SQ18T1G / NCP508SQ18T1GOSCT-ND
$ws.Range("J6").Value = "NCP508SQ18T1GOSCT-ND"
$ws.Range("I6").Value = "NCP508SQ18T1G"

# Row 4 PN filled in last
$ws.Range("I4").Value = "MAX30100EFD+"

# Column sizing for the new columns (target OOXML widths: I=14.81640625 bestFit, J=21.54296875)
$ws.Columns.Item(9).ColumnWidth = 13.98
$ws.Columns.Item(10).ColumnWidth = 20.7

# Update selection and make this sheet the active tab
$ws.Range("D9").Select() | Out-Null
$ws.Activate() | Out-Null
